$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 0.5021840000000001
$ws.Range("H2").Value = 1.506552
$ws.Range("I2").Value = 0.02664368674210811
$ws.Range("J2").Value = 0.02664368674210811
$ws.Range("M2").Value = 115.2213693333333
$ws.Range("N2").Value = 345.664108
$ws.Range("O2").Value = 0.2787408744545015
$ws.Range("P2").Value = 0.2787408744545015
$ws.Range("Q2").Value = 57.86232813729067
$ws.Range("R2").Value = 520.7609532356161
$ws.Range("S2").Value = 0.007426684541187024
$ws.Range("T2").Value = 0.007426684541187022

$ws.Range("G3").Value = 0.5021840000000001
$ws.Range("H3").Value = 1.506552
$ws.Range("I3").Value = 0.02664368674210811
$ws.Range("J3").Value = 0.02664368674210811
$ws.Range("O3").Value = 0.44716501655323
$ws.Range("P3").Value = 0.4471650165532299
$ws.Range("Q3").Value = 92.82459549556802
$ws.Range("R3").Value = 835.4213594601121
$ws.Range("S3").Value = 0.01191412462307385
$ws.Range("T3").Value = 0.01191412462307385

$ws.Range("G4").Value = 0.5021840000000001
$ws.Range("H4").Value = 1.506552
$ws.Range("I4").Value = 0.02664368674210811
$ws.Range("J4").Value = 0.02664368674210811
$ws.Range("M4").Value = 60.55095666666667
$ws.Range("N4").Value = 181.65287
$ws.Range("O4").Value = 0.1464834753134679
$ws.Range("P4").Value = 0.1464834753134678
$ws.Range("Q4").Value = 30.40772162269334
$ws.Range("R4").Value = 273.66949460424
$ws.Range("S4").Value = 0.003902859829147365
$ws.Range("T4").Value = 0.003902859829147364

$ws.Range("G5").Value = 0.5021840000000001
$ws.Range("H5").Value = 1.506552
$ws.Range("I5").Value = 0.02664368674210811
$ws.Range("J5").Value = 0.02664368674210811
$ws.Range("M5").Value = 52.74960833333333
$ws.Range("N5").Value = 158.248825
$ws.Range("O5").Value = 0.1276106336788006
$ws.Range("P5").Value = 0.1276106336788006
$ws.Range("Q5").Value = 26.49000931126667
$ws.Range("R5").Value = 238.4100838014
$ws.Range("S5").Value = 0.003400017748699876
$ws.Range("T5").Value = 0.003400017748699875

$ws.Range("I6").Value = 0.6336514295210738
$ws.Range("J6").Value = 0.6336514295210738
$ws.Range("M6").Value = 115.2213693333333
$ws.Range("N6").Value = 345.664108
$ws.Range("O6").Value = 0.2787408744545015
$ws.Range("P6").Value = 0.2787408744545015
$ws.Range("Q6").Value = 1376.106366003262
$ws.Range("R6").Value = 12384.95729402936
$ws.Range("S6").Value = 0.176624553564049
$ws.Range("T6").Value = 0.176624553564049

$ws.Range("I7").Value = 0.6336514295210738
$ws.Range("J7").Value = 0.6336514295210738
$ws.Range("O7").Value = 0.44716501655323
$ws.Range("P7").Value = 0.4471650165532299
$ws.Range("S7").Value = 0.2833467519707688
$ws.Range("T7").Value = 0.2833467519707688

$ws.Range("I8").Value = 0.6336514295210738
$ws.Range("J8").Value = 0.6336514295210738
$ws.Range("M8").Value = 60.55095666666667
$ws.Range("N8").Value = 181.65287
$ws.Range("O8").Value = 0.1464834753134679
$ws.Range("P8").Value = 0.1464834753134678
$ws.Range("Q8").Value = 723.1692993990657
$ws.Range("R8").Value = 6508.523694591591
$ws.Range("S8").Value = 0.09281946353359384
$ws.Range("T8").Value = 0.09281946353359381

$ws.Range("I9").Value = 0.6336514295210738
$ws.Range("J9").Value = 0.6336514295210738
$ws.Range("M9").Value = 52.74960833333333
$ws.Range("N9").Value = 158.248825
$ws.Range("O9").Value = 0.1276106336788006
$ws.Range("P9").Value = 0.1276106336788006
$ws.Range("Q9").Value = 629.9966078486696
$ws.Range("R9").Value = 5669.969470638026
$ws.Range("S9").Value = 0.0808606604526621
$ws.Range("T9").Value = 0.0808606604526621

$ws.Range("G10").Value = 5.487855333333333
$ws.Range("H10").Value = 16.463566
$ws.Range("I10").Value = 0.2911616028932436
$ws.Range("J10").Value = 0.2911616028932436
$ws.Range("M10").Value = 115.2213693333333
$ws.Range("N10").Value = 345.664108
$ws.Range("O10").Value = 0.2787408744545015
$ws.Range("P10").Value = 0.2787408744545015
$ws.Range("Q10").Value = 632.3182062099031
$ws.Range("R10").Value = 5690.863855889128
$ws.Range("S10").Value = 0.08115863979803703
$ws.Range("T10").Value = 0.08115863979803702

$ws.Range("G11").Value = 5.487855333333333
$ws.Range("H11").Value = 16.463566
$ws.Range("I11").Value = 0.2911616028932436
$ws.Range("J11").Value = 0.2911616028932436
$ws.Range("O11").Value = 0.44716501655323
$ws.Range("P11").Value = 0.4471650165532299
$ws.Range("Q11").Value = 1014.385068928644
$ws.Range("R11").Value = 9129.465620357796
$ws.Range("S11").Value = 0.1301972829774222
$ws.Range("T11").Value = 0.1301972829774222

$ws.Range("G12").Value = 5.487855333333333
$ws.Range("H12").Value = 16.463566
$ws.Range("I12").Value = 0.2911616028932436
$ws.Range("J12").Value = 0.2911616028932436
$ws.Range("M12").Value = 60.55095666666667
$ws.Range("N12").Value = 181.65287
$ws.Range("O12").Value = 0.1464834753134679
$ws.Range("P12").Value = 0.1464834753134678
$ws.Range("Q12").Value = 332.2948904816022
$ws.Range("R12").Value = 2990.65401433442
$ws.Range("S12").Value = 0.04265036346964218
$ws.Range("T12").Value = 0.04265036346964217

$ws.Range("G13").Value = 5.487855333333333
$ws.Range("H13").Value = 16.463566
$ws.Range("I13").Value = 0.2911616028932436
$ws.Range("J13").Value = 0.2911616028932436
$ws.Range("M13").Value = 52.74960833333333
$ws.Range("N13").Value = 158.248825
$ws.Range("O13").Value = 0.1276106336788006
$ws.Range("P13").Value = 0.1276106336788006
$ws.Range("Q13").Value = 289.4822194233278
$ws.Range("R13").Value = 2605.33997480995
$ws.Range("S13").Value = 0.03715531664814212
$ws.Range("T13").Value = 0.03715531664814212

$ws.Range("G14").Value = 0.9149506666666666
$ws.Range("H14").Value = 2.744852
$ws.Range("I14").Value = 0.04854328084357454
$ws.Range("J14").Value = 0.04854328084357455
$ws.Range("M14").Value = 115.2213693333333
$ws.Range("N14").Value = 345.664108
$ws.Range("O14").Value = 0.2787408744545015
$ws.Range("P14").Value = 0.2787408744545015
$ws.Range("Q14").Value = 105.4218686857795
$ws.Range("R14").Value = 948.7968181720159
$ws.Range("S14").Value = 0.01353099655122842
$ws.Range("T14").Value = 0.01353099655122842

$ws.Range("G15").Value = 0.9149506666666666
$ws.Range("H15").Value = 2.744852
$ws.Range("I15").Value = 0.04854328084357454
$ws.Range("J15").Value = 0.04854328084357455
$ws.Range("O15").Value = 0.44716501655323
$ws.Range("P15").Value = 0.4471650165532299
$ws.Range("Q15").Value = 169.121129967768
$ws.Range("R15").Value = 1522.090169709912
$ws.Range("S15").Value = 0.0217068569819651
$ws.Range("T15").Value = 0.0217068569819651

$ws.Range("G16").Value = 0.9149506666666666
$ws.Range("H16").Value = 2.744852
$ws.Range("I16").Value = 0.04854328084357454
$ws.Range("J16").Value = 0.04854328084357455
$ws.Range("M16").Value = 60.55095666666667
$ws.Range("N16").Value = 181.65287
$ws.Range("O16").Value = 0.1464834753134679
$ws.Range("P16").Value = 0.1464834753134678
$ws.Range("Q16").Value = 55.40113816947111
$ws.Range("R16").Value = 498.61024352524
$ws.Range("S16").Value = 0.007110788481084489
$ws.Range("T16").Value = 0.007110788481084489

$ws.Range("G17").Value = 0.9149506666666666
$ws.Range("H17").Value = 2.744852
$ws.Range("I17").Value = 0.04854328084357454
$ws.Range("J17").Value = 0.04854328084357455
$ws.Range("M17").Value = 52.74960833333333
$ws.Range("N17").Value = 158.248825
$ws.Range("O17").Value = 0.1276106336788006
$ws.Range("P17").Value = 0.1276106336788006
$ws.Range("Q17").Value = 48.26328931098889
$ws.Range("R17").Value = 434.3696037989
$ws.Range("S17").Value = 0.006194638829296531
$ws.Range("T17").Value = 0.006194638829296532
